$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the SQL text in B7: ORDER BY prt.participant_id ASC -> ORDER BY srv.survival_id ASC
$old = $ws.Range("B7").Value
$new = $old -replace "ORDER BY`n    prt.participant_id ASC", "ORDER BY`n    srv.survival_id ASC"
$ws.Range("B7").Value = $new

Write-Host $new
